# MalwareSyllabus.docx edits
# "new stuff for 4440"

$d = $word.ActiveDocument

# 1. Instructor "William L. Harrison, Ph.D" -- merge the proofErr-split
#    "Ph.D" run back into the preceding run (no text change, just clean-up).
$d.Content.Find.Execute("William L. Harrison, Ph.D", $true, $false, $false, $false, $false, `
    $true, 1, $false, "William L. Harrison, Ph.D", 2) | Out-Null

# 2. Office Hours "MW 4-5" -> "By appointment only."
$d.Content.Find.Execute("MW 4-5", $true, $false, $false, $false, $false, `
    $true, 1, $false, "By appointment only.", 2) | Out-Null

# 3. "Textbooks:" -> "Textbook:"
$d.Content.Find.Execute("Textbooks:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Textbook:", 2) | Out-Null

# 4. ", Peter Szor, Symantec Press/Addison Wesley, 2005." -- merge the
#    proofErr-split "Szor" run back in (no text change).
$d.Content.Find.Execute(", Peter Szor, Symantec Press/Addison Wesley, 2005.", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", Peter Szor, Symantec Press/Addison Wesley, 2005.", 2) | Out-Null

# 5. Prerequisites "CS3280, ECE 3210 or equivalent." -- remove gramStart/gramEnd
#    proofErr wrapper (no text change).
$d.Content.Find.Execute("CS3280, ECE 3210 or equivalent.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "CS3280, ECE 3210 or equivalent.", 2) | Out-Null

# 6. Undergraduate evaluation paragraph:
#    "...a final exam (30%), various pop quizzes and class participation (5%),
#     and a programming assignment (25% total)."
# -> "...a final exam (35%), various pop quizzes and class participation (10%),
#     and programming assignments (15% total)."
$d.Content.Find.Execute("a final exam (30%), various pop quizzes and class participation (5%), and a programming assignment (25% total).", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "a final exam (35%), various pop quizzes and class participation (10%), and programming assignments (15% total).", 2) | Out-Null

# 7. Graduate evaluation paragraph:
#    "...one final exam (30%), a programming assignment (20%), and the
#     presentation (10%)."
# -> "...one final exam (35%), programming assignments (15%), and the pop
#     quizzes and presentation (10%)."
$d.Content.Find.Execute("one final exam (30%), a programming assignment (20%), and the presentation (10%).", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "one final exam (35%), programming assignments (15%), and the pop quizzes and presentation (10%).", 2) | Out-Null

# 8. Schedule table: "using patterns; regular expressions and lex." -- merge
#    the proofErr-split "lex" run back in (no text change).
$d.Content.Find.Execute("using patterns; regular expressions and lex.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "using patterns; regular expressions and lex.", 2) | Out-Null

# 9. Schedule table: "Exam 2; Encrypted and oligomorphic viruses." -- merge
#    the proofErr-split "oligomorphic" run back in (no text change).
$d.Content.Find.Execute("Exam 2; Encrypted and oligomorphic viruses.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Exam 2; Encrypted and oligomorphic viruses.", 2) | Out-Null

$d.Saved = $false
Write-Output "edits applied"
